# Updates the cryptos list (Coin price / 1h-volume table) on Sheet1.
# Columns: A=rank (unchanged), B=Coin, C=Link, D=Price, E=Volume(1h)
# Each row in $updates is (RowNumber, Coin, Link, Price, Volume); $null means
# "leave this field unchanged" for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @(2, $null, $null, '30.113.88', '  -0.09%  '),
    @(3, $null, $null, '1.876.60', '  -1.03%  '),
    @(4, $null, $null, '0.9999', '  +0.30%  '),
    @(5, $null, $null, '242.24', '  -2.01%  '),
    @(6, $null, $null, '1.000', '  +0.38%  '),
    @(7, $null, $null, '0.4893', '  -1.86%  '),
    @(8, $null, $null, '0.2906', '  -1.60%  '),
    @(9, $null, $null, '0.06585', '  -1.23%  '),
    @(10, $null, $null, '1.879.11', '  -0.69%  '),
    @(11, $null, $null, '16.39', '  -3.83%  '),
    @(12, $null, $null, '0.07200', '  +0.03%  '),
    @(13, $null, $null, '0.6655', '  -2.17%  '),
    @(14, $null, $null, '4.923', '  +1.32%  '),
    @(15, $null, $null, '85.87', '  -0.15%  '),
    @(16, $null, $null, '30.084.78', '  -0.14%  '),
    @(17, $null, $null, '0.000007786', '  -3.18%  '),
    @(18, $null, $null, '1.000', '  +0.22%  '),
    @(19, $null, $null, '12.76', '  -1.57%  '),
    @(20, $null, $null, '2.119.49', '  -0.73%  '),
    @(21, $null, $null, '0.9954', '  -0.13%  '),
    @(22, $null, $null, '4.770', '  -0.22%  '),
    @(23, $null, $null, '5.819', '  +2.39%  '),
    @(24, $null, $null, '9.197', '  +0.05%  '),
    @(25, $null, $null, '152.87', '  +3.91%  '),
    @(26, $null, $null, '143.11', '  +7.20%  '),
    @(27, $null, $null, $null, '  +0.44%  '),
    @(28, $null, $null, '1.883', '  -3.58%  '),
    @(29, $null, $null, '1.398', '  +2.01%  '),
    @(30, $null, $null, '4.206', '  -0.81%  '),
    @(31, $null, $null, '0.08781', '  +0.30%  '),
    @(32, $null, $null, '4.008', '  +1.13%  '),
    @(33, $null, $null, '0.05127', '  -0.09%  '),
    @(34, $null, $null, '0.7160', '  +0.90%  '),
    @(35, $null, $null, '1.106', '  -1.41%  '),
    @(36, $null, $null, '2.665', '  +0.23%  '),
    @(37, $null, $null, '0.01831', '  +9.84%  '),
    @(38, $null, $null, '2.676', '  -4.40%  '),
    @(39, $null, $null, '2.136', '  -4.70%  '),
    @(40, $null, $null, '0.9244', '  -1.55%  '),
    @(41, $null, $null, '5.800', '  -4.72%  '),
    @(42, $null, $null, '0.9997', '  +0.42%  '),
    @(43, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '103.76', '  +0.57%  '),
    @(44, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.4217', '  -0.47%  '),
    @(45, $null, $null, '7.408', '  -1.38%  '),
    @(46, $null, $null, '0.1280', '  +1.15%  '),
    @(47, $null, $null, '0.05693', '  -0.50%  '),
    @(48, $null, $null, '32.80', '  -0.13%  '),
    @(49, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '8.235', '  -0.58%  '),
    @(50, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.3748', '  -0.18%  '),
    @(51, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '55.73', '  -0.71%  ')
)

foreach ($u in $updates) {
    $row    = $u[0]
    $coin   = $u[1]
    $link   = $u[2]
    $price  = $u[3]
    $volume = $u[4]

    if ($null -ne $coin) {
        $ws.Cells.Item($row, 2).Value = $coin
    }
    if ($null -ne $link) {
        $ws.Cells.Item($row, 3).Value = $link
    }
    if ($null -ne $price) {
        # Force the cell to text so Excel doesn't auto-convert numeric-looking
        # price strings (e.g. "0.9999", "242.24") into real numbers, then
        # restore the default (unstyled) cell style so formatting matches
        # the rest of the sheet.
        $priceCell = $ws.Cells.Item($row, 4)
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $price
        $priceCell.Style = "Normal"
    }
    if ($null -ne $volume) {
        $ws.Cells.Item($row, 5).Value = $volume
    }
}

